# Generate Report for handback
# Update the "Correspond Handoff Datetime" (D2) and "Correspond Handback DateTime" (G2)
# values on the zh-cn and de-de worksheets to reflect the newly generated report times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-19 05:48:23"
$wsZhCn.Range("G2").Value = "2016-02-19 05:49:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-19 05:48:33"
$wsDeDe.Range("G2").Value = "2016-02-19 05:49:30"
